# feat: add 2022-Q1 data
#
# The existing "总计" (Totals) sheet holds the running per-quarter summary
# table in columns A:D. We need to:
#   1. Turn the *current* "总计" sheet into the new "2022-Q1" holdings sheet
#      (reusing its sheetId/rId, matching how Excel keeps the original sheet
#      identity when it is renamed) and replace its contents with the
#      2022-Q1 per-fund holdings table (columns A:H).
#   2. Insert a brand-new sheet named "总计" right after it, containing the
#      updated summary table (the new 2022-Q1 row on top of the previous
#      quarters' rows).

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")

# Grab a style donor cell (style index 2: bold, centered, thin-bordered)
# from an existing, already-styled sheet so the new sheets reuse the same
# cell style instead of minting new ones.
$donor = $q4.Cells.Item(1, 2)

# --- Step 1: rename the existing "总计" sheet to "2022-Q1" -----------------
$totals.Name = "2022-Q1"
$q1 = $totals

# --- Step 2: insert the new "总计" sheet right after it --------------------
$newTotals = $wb.Worksheets.Add($null, $q1)
$newTotals.Name = "总计"

$headers1 = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$rows1 = @(
    @("002906", "南方中证500量化增强股票A", "7.21", "92.26", "1.13", "0.0815", 4),
    @("001364", "大成景润灵活配置混合", "4.72", "26.73", "1.42", "0.0670", 10),
    @("008114", "天弘中证红利低波动100指数A", "3.16", "92.60", "2.08", "0.0657", 3),
    @("008115", "天弘中证红利低波动100指数C", "2.37", "92.60", "2.08", "0.0493", 3),
    @("515100", "景顺长城中证红利低波动100ETF", "1.25", "97.96", "2.21", "0.0276", 3),
    @("002907", "南方中证500量化增强股票C", "1.36", "92.26", "1.13", "0.0154", 4)
)

$headers2 = @("日期", "持有数量(只)", "持有市值(亿元)")
$rows2 = @(
    @("2022-Q1", 6, 0.31),
    @("2021-Q4", 3, 0.06),
    @("2021-Q3", 3, 0.25),
    @("2021-Q2", 3, 0.31),
    @("2021-Q1", 8, 1.09),
    @("2020-Q4", 12, 0.72)
)

# ============================================================================
# Phase A: clone the "index/header" cell style (s=2) onto every cell that
# needs it FIRST. (ClearFormats(), used in phase B below to drop the stray
# text-NumberFormat style off data cells, invalidates the pending
# Copy()/PasteSpecial clipboard state, so all format-paste work must finish
# before any ClearFormats() call happens.)
# ============================================================================

$donor.Copy()

for ($i = 0; $i -lt $headers1.Length; $i++) {
    $q1.Cells.Item(1, 2 + $i).PasteSpecial(-4122)
}
for ($r = 0; $r -lt $rows1.Length; $r++) {
    $q1.Cells.Item($r + 2, 1).PasteSpecial(-4122)
}

for ($i = 0; $i -lt $headers2.Length; $i++) {
    $newTotals.Cells.Item(1, 2 + $i).PasteSpecial(-4122)
}
for ($r = 0; $r -lt $rows2.Length; $r++) {
    $newTotals.Cells.Item($r + 2, 1).PasteSpecial(-4122)
}

# ============================================================================
# Phase B: write the actual values.
# ============================================================================

# --- "2022-Q1": per-fund holdings table -------------------------------------
for ($i = 0; $i -lt $headers1.Length; $i++) {
    $q1.Cells.Item(1, 2 + $i).Value = $headers1[$i]
}

for ($r = 0; $r -lt $rows1.Length; $r++) {
    $excelRow = $r + 2
    $q1.Cells.Item($excelRow, 1).Value = $r

    $data = $rows1[$r]
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $q1.Cells.Item($excelRow, 2 + $c)
        # Force text storage for numeric-looking strings ("7.21", "0.0815",
        # ...) without leaving a stray NumberFormat-derived style behind.
        $cell.NumberFormat = "@"
        $cell.Value = $data[$c]
        $cell.ClearFormats()
    }
    $q1.Cells.Item($excelRow, 8).Value = $data[6]
}

# --- "总计": refreshed summary table ----------------------------------------
for ($i = 0; $i -lt $headers2.Length; $i++) {
    $newTotals.Cells.Item(1, 2 + $i).Value = $headers2[$i]
}

for ($r = 0; $r -lt $rows2.Length; $r++) {
    $excelRow = $r + 2
    $newTotals.Cells.Item($excelRow, 1).Value = $r

    $data = $rows2[$r]
    $newTotals.Cells.Item($excelRow, 2).Value = $data[0]
    $newTotals.Cells.Item($excelRow, 3).Value = $data[1]
    $newTotals.Cells.Item($excelRow, 4).Value = $data[2]
}
